$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Treatment query text in cell B5: drop the redundant CONCAT() wrapper
# around REPLACE() for the "Treatment Agent" column.
$b5 = $ws.Range("B5").Value2
$oldFragment = "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"","
$newFragment = "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","
$ws.Range("B5").Value2 = $b5.Replace($oldFragment, $newFragment)

# Match the author's final formatting pass on B5 (font size bump + wrap).
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true

# Restore the view to the top of the sheet with C5 selected (matches the saved
# selection/scroll position in the workbook).
$ws.Range("C5").Select() | Out-Null
